# Auto-generated edit script applying numeric updates from the commit diff
# to the Halicarnassus Profits workbook. Each worksheet (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR) gets specific cell values updated or cleared to match
# the refreshed market-board pricing data pulled in by the scheduled runner.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1500
$ws.Range("I18").Value = 1500
$ws.Range("K18").Value = 1500
$ws.Range("M18").Value = -1216
$ws.Range("H32").Value = 700
$ws.Range("I32").Value = 700
$ws.Range("K32").Value = 700
$ws.Range("M32").Value = -374
$ws.Range("H64").Value = 6314.154
$ws.Range("J64").Value = 8000
$ws.Range("L64").Value = 8000
$ws.Range("N64").Value = -8496
$ws.Range("H67").Value = 6314.154
$ws.Range("J67").Value = 8000
$ws.Range("L67").Value = 8000
$ws.Range("N67").Value = -9716
$ws.Range("H141").Value = 1442.6666
$ws.Range("I141").Value = 1442.6666
$ws.Range("K141").Value = 4327.9998
$ws.Range("M141").Value = 852.0002000000004

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 2376.6667
$ws.Range("I19").Value = 2376.6667
$ws.Range("K19").Value = 2376.6667
$ws.Range("M19").Value = -2147.6667
$ws.Range("H125").Value = 67500
$ws.Range("J125").Value = 67500
$ws.Range("L125").Value = 67500
$ws.Range("N125").Value = -77340

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2089.4211
$ws.Range("I20").Value = 2150.611
$ws.Range("K20").Value = 2150.611
$ws.Range("M20").Value = -1903.611

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 10000000
$ws.Range("I3").Value = 10000000
$ws.Range("K3").Value = 10000000
$ws.Range("M3").Value = -9999887
$ws.Range("H12").Value = 166
$ws.Range("I12").Value = 112.4
$ws.Range("J12").Value = 300
$ws.Range("K12").Value = 112.4
$ws.Range("L12").Value = 300
$ws.Range("M12").Value = 57.59999999999999
$ws.Range("N12").Value = -640
$ws.Range("H31").Value = 8086.76
$ws.Range("I31").Value = 4703.2
$ws.Range("J31").Value = 8932.65
$ws.Range("K31").Value = 4703.2
$ws.Range("L31").Value = 8932.65
$ws.Range("M31").Value = -4408.2
$ws.Range("N31").Value = -9522.65
$ws.Range("H34").Value = 8086.76
$ws.Range("I34").Value = 4703.2
$ws.Range("J34").Value = 8932.65
$ws.Range("K34").Value = 4703.2
$ws.Range("L34").Value = 8932.65
$ws.Range("M34").Value = -4501.2
$ws.Range("N34").Value = -9336.65
$ws.Range("H132").Value = 4027.08
$ws.Range("I132").Value = 3565.6
$ws.Range("J132").Value = 4719.3
$ws.Range("K132").Value = 10696.8
$ws.Range("L132").Value = 14157.9
$ws.Range("M132").Value = -8166.799999999999
$ws.Range("N132").Value = -19217.9
$ws.Range("H134").Value = 854.2857
$ws.Range("I134").Value = 854.2857
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 2562.8571
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -27.85710000000017
$ws.Range("N134").ClearContents()
$ws.Range("H141").Value = 73537.5
$ws.Range("J141").Value = 73537.5
$ws.Range("L141").Value = 73537.5
$ws.Range("N141").Value = -83897.5

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 376
$ws.Range("I14").Value = 376
$ws.Range("K14").Value = 1128
$ws.Range("M14").Value = -955
$ws.Range("H60").Value = 1313.2667
$ws.Range("I60").Value = 300
$ws.Range("J60").Value = 2833.1667
$ws.Range("K60").Value = 900
$ws.Range("L60").Value = 8499.500100000001
$ws.Range("M60").Value = -649
$ws.Range("N60").Value = -9001.500100000001
$ws.Range("H81").Value = 2225
$ws.Range("I81").Value = 950
$ws.Range("K81").Value = 2850
$ws.Range("M81").Value = -1727
$ws.Range("H84").Value = 2225
$ws.Range("I84").Value = 950
$ws.Range("K84").Value = 8550
$ws.Range("M84").Value = -2934
$ws.Range("H101").Value = 5000
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 40.714287
$ws.Range("I2").Value = 19.333334
$ws.Range("K2").Value = 19.333334
$ws.Range("M2").Value = 93.66666599999999
$ws.Range("H132").Value = 4172.4517
$ws.Range("I132").Value = 3944.5
$ws.Range("J132").Value = 5357.8
$ws.Range("K132").Value = 11833.5
$ws.Range("L132").Value = 16073.4
$ws.Range("M132").Value = -9303.5
$ws.Range("N132").Value = -21133.4

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 35980
$ws.Range("J5").Value = 35980
$ws.Range("L5").Value = 35980
$ws.Range("N5").Value = -36206
$ws.Range("H14").Value = 18000
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H22").Value = 783.6923
$ws.Range("I22").Value = 769.3
$ws.Range("J22").Value = 831.6667
$ws.Range("K22").Value = 769.3
$ws.Range("L22").Value = 831.6667
$ws.Range("M22").Value = -474.3
$ws.Range("N22").Value = -1421.6667
$ws.Range("H26").Value = 60000
$ws.Range("I26").Value = 60000
$ws.Range("K26").Value = 60000
$ws.Range("M26").Value = -59705
$ws.Range("H27").Value = 783.6923
$ws.Range("I27").Value = 769.3
$ws.Range("J27").Value = 831.6667
$ws.Range("K27").Value = 769.3
$ws.Range("L27").Value = 831.6667
$ws.Range("M27").Value = -662.3
$ws.Range("N27").Value = -1045.6667
$ws.Range("H31").Value = 1260.75
$ws.Range("I31").Value = 674
$ws.Range("K31").Value = 674
$ws.Range("M31").Value = -426
$ws.Range("H68").Value = 4253.1763
$ws.Range("I68").Value = 2590.4
$ws.Range("K68").Value = 2590.4
$ws.Range("M68").Value = -1841.4
$ws.Range("H71").Value = 4253.1763
$ws.Range("I71").Value = 2590.4
$ws.Range("K71").Value = 12952
$ws.Range("M71").Value = -9208
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H99").Value = 20000
$ws.Range("I99").Value = 20000
$ws.Range("K99").Value = 20000
$ws.Range("M99").Value = -17005
$ws.Range("H100").Value = 7488.75
$ws.Range("I100").Value = 2985
$ws.Range("J100").Value = 8990
$ws.Range("K100").Value = 2985
$ws.Range("L100").Value = 8990
$ws.Range("M100").Value = -2444
$ws.Range("N100").Value = -10072
$ws.Range("H122").Value = 1333.3334
$ws.Range("I122").Value = 1333.3334
$ws.Range("K122").Value = 4000.0002
$ws.Range("M122").Value = -1550.0002
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 128138
$ws.Range("I11").Value = 749.5
$ws.Range("J11").Value = 255526.5
$ws.Range("K11").Value = 749.5
$ws.Range("L11").Value = 255526.5
$ws.Range("M11").Value = -607.5
$ws.Range("N11").Value = -255810.5
$ws.Range("H17").Value = 50502.5
$ws.Range("I17").Value = 100000
$ws.Range("J17").Value = 1005
$ws.Range("K17").Value = 100000
$ws.Range("L17").Value = 1005
$ws.Range("M17").Value = -99828
$ws.Range("N17").Value = -1349
$ws.Range("H106").Value = 68999
$ws.Range("J106").Value = 68999
$ws.Range("L106").Value = 68999
$ws.Range("N106").Value = -71523
$ws.Range("H122").Value = 1431.7858
$ws.Range("I122").Value = 1462.7273
$ws.Range("J122").Value = 1318.3334
$ws.Range("K122").Value = 4388.1819
$ws.Range("L122").Value = 3955.0002
$ws.Range("M122").Value = -1938.1819
$ws.Range("N122").Value = -8855.0002
